$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6142318033764959
$ws.Range("J2").Value = 0.6142318033764957
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 18.15977915022644
$ws.Range("R2").Value = 163.438012352038
$ws.Range("S2").Value = 0.1412704572241715
$ws.Range("T2").Value = 0.1412704572241716

# Row 3
$ws.Range("I3").Value = 0.6142318033764959
$ws.Range("J3").Value = 0.6142318033764957
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("S3").Value = 0.4173540458919317
$ws.Range("T3").Value = 0.4173540458919317

# Row 4
$ws.Range("I4").Value = 0.6142318033764959
$ws.Range("J4").Value = 0.6142318033764957
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 7.148106629729777
$ws.Range("R4").Value = 64.33295966756799
$ws.Range("S4").Value = 0.05560730026039252
$ws.Range("T4").Value = 0.05560730026039252

# Row 5
$ws.Range("G5").Value = 1.102210333333334
$ws.Range("H5").Value = 3.306631
$ws.Range("I5").Value = 0.3857681966235041
$ws.Range("J5").Value = 0.3857681966235041
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 11.40524670874122
$ws.Range("R5").Value = 102.647220378671
$ws.Range("S5").Value = 0.08872489053801401
$ws.Range("T5").Value = 0.08872489053801402

# Row 6
$ws.Range("G6").Value = 1.102210333333334
$ws.Range("H6").Value = 3.306631
$ws.Range("I6").Value = 0.3857681966235041
$ws.Range("J6").Value = 0.3857681966235041
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("Q6").Value = 33.69441815237745
$ws.Range("R6").Value = 303.249763371397
$ws.Range("S6").Value = 0.2621191490772856
$ws.Range("T6").Value = 0.2621191490772856

# Row 7
$ws.Range("G7").Value = 1.102210333333334
$ws.Range("H7").Value = 3.306631
$ws.Range("I7").Value = 0.3857681966235041
$ws.Range("J7").Value = 0.3857681966235041
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 4.48936735067289
$ws.Range("R7").Value = 40.404306156056
$ws.Range("S7").Value = 0.0349241570082045
$ws.Range("T7").Value = 0.0349241570082045
